# Updates cryptos list data (prices in column D, 1h volume % in column E)
# Also reflects a rank swap between Aave and BabyDogeCoin at rows 46/47.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.616.69"
$ws.Range("D3").Value = "1.697.58"
$ws.Range("E3").Value = "  -5.99%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'219.85"
$ws.Range("E5").Value = "  -5.45%  "
$ws.Range("D6").Value = "'0.5142"
$ws.Range("E6").Value = "  -13.22%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.2648"
$ws.Range("E8").Value = "  -4.58%  "
$ws.Range("D9").Value = "'22.18"
$ws.Range("D10").Value = "'0.06256"
$ws.Range("E10").Value = "  -8.35%  "
$ws.Range("D11").Value = "'0.07343"
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("D12").Value = "1.696.73"
$ws.Range("E12").Value = "  -6.24%  "
$ws.Range("D13").Value = "'4.520"
$ws.Range("E13").Value = "  -4.92%  "
$ws.Range("D14").Value = "'0.5849"
$ws.Range("E14").Value = "  -6.08%  "
$ws.Range("D15").Value = "1.928.64"
$ws.Range("E15").Value = "  -5.96%  "
$ws.Range("D16").Value = "'0.000008401"
$ws.Range("E16").Value = "  -9.41%  "
$ws.Range("D17").Value = "'65.65"
$ws.Range("E17").Value = "  -13.15%  "
$ws.Range("D18").Value = "26.660.22"
$ws.Range("E18").Value = "  -7.07%  "
$ws.Range("D19").Value = "'5.036"
$ws.Range("E19").Value = "  -7.96%  "
$ws.Range("D20").Value = "'1.005"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  -4.97%  "
$ws.Range("D22").Value = "'186.92"
$ws.Range("E22").Value = "  -11.37%  "
$ws.Range("D23").Value = "'6.280"
$ws.Range("E23").Value = "  -7.94%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "'145.05"
$ws.Range("E25").Value = "  -5.79%  "
$ws.Range("D26").Value = "'7.608"
$ws.Range("E26").Value = "  -3.33%  "
$ws.Range("D27").Value = "'0.1151"
$ws.Range("E27").Value = "  -9.26%  "
$ws.Range("D28").Value = "'15.79"
$ws.Range("E28").Value = "  -3.94%  "
$ws.Range("D29").Value = "'1.326"
$ws.Range("E29").Value = "  -7.52%  "
$ws.Range("D30").Value = "'0.05683"
$ws.Range("E30").Value = "  -7.64%  "
$ws.Range("D31").Value = "'1.338"
$ws.Range("E31").Value = "  -6.35%  "
$ws.Range("D32").Value = "'3.516"
$ws.Range("E32").Value = "  -7.01%  "
$ws.Range("D33").Value = "'3.519"
$ws.Range("E33").Value = "  -6.12%  "
$ws.Range("D34").Value = "'1.663"
$ws.Range("E34").Value = "  -4.01%  "
$ws.Range("D35").Value = "'1.026"
$ws.Range("E35").Value = "  -3.34%  "
$ws.Range("D36").Value = "'0.6026"
$ws.Range("E36").Value = "  -6.15%  "
$ws.Range("D37").Value = "'2.374"
$ws.Range("E37").Value = "  -4.97%  "
$ws.Range("D38").Value = "'2.683"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").Value = "1.100.72"
$ws.Range("E39").Value = "  -3.91%  "
$ws.Range("E40").Value = "  -5.79%  "
$ws.Range("D41").Value = "'0.8624"
$ws.Range("E41").Value = "  -2.32%  "
$ws.Range("D42").Value = "'5.876"
$ws.Range("E42").Value = "  -10.60%  "
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "'99.01"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("D45").Value = "1.857.87"
$ws.Range("E45").Value = "  -5.27%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000107"
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'56.76"
$ws.Range("E47").Value = "  -6.15%  "
$ws.Range("D48").Value = "'8.194"
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("D49").Value = "'1.004"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "'0.05244"
$ws.Range("E50").Value = "  -4.08%  "
$ws.Range("D51").Value = "'0.4326"
$ws.Range("E51").Value = "  -3.41%  "
